# Insert a new row at row 9 on the active sheet, shifting existing rows 9-19
# down to 10-20, then populate the new row 9 with the new price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 9 (pushes old rows 9..19 -> 10..20)
$ws.Rows("9:9").Insert()

# Populate the newly inserted row 9 with the new record
$ws.Range("A9").Value = 10
$ws.Range("B9").Value = "Vega Modelo de Temuco"
$ws.Range("C9").Value = "La Araucanía"
$ws.Range("D9").Value = 45093
$ws.Range("E9").Value = 9
$ws.Range("F9").Value = "Fruta"
$ws.Range("G9").Value = 100104
$ws.Range("H9").Value = "Frutos de pepita"
$ws.Range("I9").Value = 100104005
$ws.Range("J9").Value = "Pera asiática"
$ws.Range("K9").Value = "Hosui"
$ws.Range("L9").Value = "Primera"
$ws.Range("M9").Value = 170
$ws.Range("N9").Value = 15000
$ws.Range("O9").Value = 16000
$ws.Range("P9").Value = 15471
$ws.Range("Q9").Value = "$/caja 18 kilos granel"
$ws.Range("R9").Value = "Región de O'Higgins"
$ws.Range("S9").Value = 860
$ws.Range("T9").Value = 18

# Make sure the D9 cell keeps the date number format used by the rest of the
# date column (carried over from the insert, but set explicitly to be safe).
$ws.Range("D9").NumberFormat = "YYYY-MM-DD HH:MM:SS"
